# Re-activate spherical coordinates: extend the "dx"/"dtmax" helper table on
# Sheet1 from a single sample column (B) to three sample columns (B, C, D),
# and bump the budget/hours-spent numbers on the "begroting" sheet.

$wb = $excel.ActiveWorkbook

# --- "begroting" sheet: hours spent on Part 5 grows from 30 to 33 ---
$ws1 = $wb.Worksheets.Item("begroting")
$ws1.Range("B40").Value = 33

# B43/C43 already hold formulas (SUM(B27:B40) and B43*135); simply recalculate
# everything so their cached values follow the new B40 input.
$wb.Application.Calculate()

# --- "Sheet1": add two more spherical-coordinate samples (dx = 3 and dx = 1) ---
$ws2 = $wb.Worksheets.Item("Sheet1")

$ws2.Range("C6").Value = 3
$ws2.Range("D6").Value = 1

$ws2.Range("C7").Formula = "=C6/SQRT(9.81)"
$ws2.Range("D7").Formula = "=D6/SQRT(9.81)"

# Update the active selection like Excel would after clicking around L16
$ws2.Range("L16").Select()

$wb.Application.Calculate()
